$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 135; this shifts all rows 135..201 down to 136..202
$ws.Rows(135).Insert()

# Fill the new row 135 with its data (a new weekly record)
$ws.Range("A135").Value = 8
$ws.Range("B135").Value = "Terminal La Palmera de La Serena"
$ws.Range("C135").Value = "Coquimbo"
$ws.Range("D135").Value = 45126
$ws.Range("E135").Value = 4
$ws.Range("F135").Value = 100112052
$ws.Range("G135").Value = "Albahaca"
$ws.Range("H135").Value = "Sin especificar"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 3000
$ws.Range("L135").Value = 3500
$ws.Range("M135").Value = 3250
$ws.Range("N135").Value = "$/paquete"
$ws.Range("O135").Value = "Región de Arica y Parinacota"
$ws.Range("P135").Value = 3250
$ws.Range("Q135").Value = 1
$ws.Range("R135").Value = "Hortaliza"
